# Refresh crypto price/volume snapshot to match latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.951.38'
$ws.Range('E2').Value = '  +4.73%  '

$ws.Range('D3').Value = '2.618.42'
$ws.Range('E3').Value = '  +5.44%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('E5').Value = '  +3.02%  '

$ws.Range('D6').Value = "'180.70"
$ws.Range('E6').Value = '  +3.87%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  +2.06%  '

$ws.Range('D9').Value = '2.617.69'
$ws.Range('E9').Value = '  +5.46%  '

$ws.Range('E10').Value = '  +14.65%  '

$ws.Range('D12').Value = "'0.345"
$ws.Range('E12').Value = '  +3.68%  '

$ws.Range('D13').Value = "'5.04"
$ws.Range('E13').Value = '  +1.80%  '

$ws.Range('E14').Value = '  +5.21%  '

$ws.Range('D15').Value = "'26.63"
$ws.Range('E15').Value = '  +5.57%  '

$ws.Range('E16').Value = '  +7.48%  '

$ws.Range('D17').Value = '70.984.97'
$ws.Range('E17').Value = '  +4.89%  '

$ws.Range('D18').Value = '2.613.74'
$ws.Range('E18').Value = '  +6.41%  '

$ws.Range('D19').Value = "'380.15"
$ws.Range('E19').Value = '  +9.95%  '

$ws.Range('E20').Value = '  +7.41%  '

$ws.Range('D21').Value = "'11.50"
$ws.Range('E21').Value = '  +6.72%  '

$ws.Range('E22').Value = '  +2.92%  '

$ws.Range('D23').Value = "'71.94"
$ws.Range('E23').Value = '  +1.71%  '

$ws.Range('E24').Value = '  +6.47%  '

$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('D26').Value = "'1.85"
$ws.Range('E26').Value = '  +10.22%  '

$ws.Range('D27').Value = "'9.63"
$ws.Range('E27').Value = '  +9.23%  '

$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = '  +0.13%  '

$ws.Range('D30').Value = '0.0₃0945'
$ws.Range('E30').Value = '  +6.65%  '

$ws.Range('D31').Value = "'524.64"
$ws.Range('E31').Value = '  +5.39%  '

$ws.Range('D32').Value = "'8.01"
$ws.Range('E32').Value = '  +4.07%  '

$ws.Range('E33').Value = '  +6.67%  '

$ws.Range('E34').Value = '  +4.29%  '

$ws.Range('E35').Value = '  +0.01%  '

$ws.Range('D36').Value = "'164.79"
$ws.Range('E36').Value = '  +0.18%  '

$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('D38').Value = "'19.10"
$ws.Range('E38').Value = '  +4.71%  '

$ws.Range('E41').Value = '  +5.31%  '

$ws.Range('D43').Value = "'5.02"
$ws.Range('E43').Value = '  +5.75%  '

$ws.Range('E44').Value = '  +9.28%  '

$ws.Range('E45').Value = '  +2.70%  '

$ws.Range('D46').Value = "'40.13"
$ws.Range('E46').Value = '  +3.90%  '

$ws.Range('D47').Value = "'153.53"
$ws.Range('E47').Value = '  +4.23%  '

$ws.Range('E48').Value = '  +3.89%  '

$ws.Range('E49').Value = '  +7.56%  '

$ws.Range('E50').Value = '  +4.58%  '

$ws.Range('E51').Value = '  +7.25%  '

# Rows 39/40 swapped order: WhiteBITCoin now ranks above Stacks, both with refreshed price/volume
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = "'18.94"
$ws.Range('E39').Value = '  +1.67%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'1.86"
$ws.Range('E40').Value = '  +8.48%  '

